# Update computed mass-flow results with new input-file-derived values
# (per commit "added new input files").
$wb = $excel.ActiveWorkbook
$wsOutput = $wb.Worksheets.Item("Output_flows")
$wsInput  = $wb.Worksheets.Item("Input_flows")

# --- Output_flows sheet updates ---
$wsOutput.Cells.Item(2, 3).Value = [double]"7.627930459241456E-13"
$wsOutput.Cells.Item(2, 5).Value = [double]"3.90645089157181E-10"
$wsOutput.Cells.Item(2, 7).Value = [double]"5.502388722896761E-10"
$wsOutput.Cells.Item(2, 9).Value = [double]"4.063066203009719E-12"
$wsOutput.Cells.Item(2, 13).Value = [double]"1.121959070153541E-28"
$wsOutput.Cells.Item(3, 3).Value = [double]"4.776009410686622E-14"
$wsOutput.Cells.Item(3, 4).Value = [double]"4.719400747970281E-14"
$wsOutput.Cells.Item(3, 5).Value = [double]"1.772998952331601E-10"
$wsOutput.Cells.Item(3, 7).Value = [double]"3.445162546018305E-11"
$wsOutput.Cells.Item(3, 9).Value = [double]"2.543972120027284E-13"
$wsOutput.Cells.Item(3, 13).Value = [double]"7.024824237833165E-30"
$wsOutput.Cells.Item(4, 3).Value = [double]"1.872491002410639E-15"
$wsOutput.Cells.Item(4, 4).Value = [double]"1.850296906360981E-14"
$wsOutput.Cells.Item(4, 5).Value = [double]"1.347405424528664E-08"
$wsOutput.Cells.Item(4, 7).Value = [double]"1.350716741643516E-12"
$wsOutput.Cells.Item(4, 9).Value = [double]"9.973943716433702E-15"
$wsOutput.Cells.Item(4, 13).Value = [double]"2.754165464880795E-31"
$wsOutput.Cells.Item(5, 3).Value = [double]"9.257468677323573E-16"
$wsOutput.Cells.Item(5, 4).Value = [double]"9.147742569835364E-14"
$wsOutput.Cells.Item(5, 5).Value = [double]"5.682861142806018E-05"
$wsOutput.Cells.Item(5, 7).Value = [double]"6.677852075979816E-13"
$wsOutput.Cells.Item(5, 9).Value = [double]"4.931050212011886E-15"
$wsOutput.Cells.Item(5, 13).Value = [double]"1.361640749700588E-31"
$wsOutput.Cells.Item(6, 3).Value = [double]"1.361315705238656E-12"
$wsOutput.Cells.Item(6, 4).Value = [double]"1.345180422624706E-09"
$wsOutput.Cells.Item(6, 5).Value = [double]"821.9782193676954"
$wsOutput.Cells.Item(6, 7).Value = [double]"9.819817085160352E-10"
$wsOutput.Cells.Item(6, 9).Value = [double]"7.251135629953755E-12"
$wsOutput.Cells.Item(6, 13).Value = [double]"2.002299982932536E-28"
$wsOutput.Cells.Item(7, 3).Value = [double]"6.888708747556553E-07"
$wsOutput.Cells.Item(7, 6).Value = [double]"3.572386075704471E-12"
$wsOutput.Cells.Item(7, 7).Value = [double]"0.004969152974114131"
$wsOutput.Cells.Item(7, 9).Value = [double]"3.669315005443532E-05"
$wsOutput.Cells.Item(7, 10).Value = [double]"0.004410849679081186"
$wsOutput.Cells.Item(7, 13).Value = [double]"1.013230168033768E-21"
$wsOutput.Cells.Item(8, 3).Value = [double]"2.781369441358133E-07"
$wsOutput.Cells.Item(8, 4).Value = [double]"2.749318500198406E-08"
$wsOutput.Cells.Item(8, 6).Value = [double]"1.028781286322824E-11"
$wsOutput.Cells.Item(8, 7).Value = [double]"0.002006333949963742"
$wsOutput.Cells.Item(8, 9).Value = [double]"1.481514315796482E-05"
$wsOutput.Cells.Item(8, 13).Value = [double]"4.090995177334069E-22"
$wsOutput.Cells.Item(9, 3).Value = [double]"1.854018461962111E-07"
$wsOutput.Cells.Item(9, 4).Value = [double]"1.832043926290027E-07"
$wsOutput.Cells.Item(9, 6).Value = [double]"1.334089506356708E-08"
$wsOutput.Cells.Item(9, 7).Value = [double]"0.00133739162039466"
$wsOutput.Cells.Item(9, 9).Value = [double]"9.875548542039827E-06"
$wsOutput.Cells.Item(9, 13).Value = [double]"2.726995009649529E-22"
$wsOutput.Cells.Item(10, 3).Value = [double]"9.257458033827716E-08"
$wsOutput.Cells.Item(10, 4).Value = [double]"9.147732055543041E-07"
$wsOutput.Cells.Item(10, 6).Value = [double]"5.682854493220766E-05"
$wsOutput.Cells.Item(10, 7).Value = [double]"0.000667784439832047"
$wsOutput.Cells.Item(10, 9).Value = [double]"4.931044542685385E-06"
$wsOutput.Cells.Item(10, 13).Value = [double]"1.361639184195137E-22"
$wsOutput.Cells.Item(11, 3).Value = [double]"0.0001359659559556395"
$wsOutput.Cells.Item(11, 4).Value = [double]"0.01343543906760231"
$wsOutput.Cells.Item(11, 6).Value = [double]"820.9782193700295"
$wsOutput.Cells.Item(11, 7).Value = [double]"0.9807870519346413"
$wsOutput.Cells.Item(11, 9).Value = [double]"0.007242314063495054"
$wsOutput.Cells.Item(11, 13).Value = [double]"1.999864030376811E-19"
$wsOutput.Cells.Item(12, 3).Value = [double]"3.881395958277108E-14"
$wsOutput.Cells.Item(12, 5).Value = [double]"4.407377887160568E-08"
$wsOutput.Cells.Item(12, 9).Value = [double]"1.033725256335287E-12"
$wsOutput.Cells.Item(12, 10).Value = [double]"2.692645902759782E-08"
$wsOutput.Cells.Item(12, 13).Value = [double]"2.854488136405577E-29"
$wsOutput.Cells.Item(13, 3).Value = [double]"2.804018764252148E-15"
$wsOutput.Cells.Item(13, 4).Value = [double]"2.078087653636733E-14"
$wsOutput.Cells.Item(13, 5).Value = [double]"1.187483406598804E-08"
$wsOutput.Cells.Item(13, 9).Value = [double]"7.467893116300725E-14"
$wsOutput.Cells.Item(13, 10).Value = [double]"3.804524440595465E-09"
$wsOutput.Cells.Item(13, 13).Value = [double]"2.062154539978769E-30"
$wsOutput.Cells.Item(14, 3).Value = [double]"4.11401575365547E-17"
$wsOutput.Cells.Item(14, 4).Value = [double]"1.219576053274894E-15"
$wsOutput.Cells.Item(14, 5).Value = [double]"1.494333837079457E-08"
$wsOutput.Cells.Item(14, 9).Value = [double]"1.095678471155684E-15"
$wsOutput.Cells.Item(14, 10).Value = [double]"3.237892007881076E-10"
$wsOutput.Cells.Item(14, 13).Value = [double]"3.025563299398062E-32"
$wsOutput.Cells.Item(15, 3).Value = [double]"1.742391472607643E-11"
$wsOutput.Cells.Item(15, 4).Value = [double]"4.390435492388066E-09"
$wsOutput.Cells.Item(15, 5).Value = [double]"5.772246289752951"
$wsOutput.Cells.Item(15, 9).Value = [double]"4.640480103084508E-10"
$wsOutput.Cells.Item(15, 13).Value = [double]"1.281403866288479E-26"
$wsOutput.Cells.Item(16, 3).Value = [double]"3.84794750988665E-11"
$wsOutput.Cells.Item(16, 4).Value = [double]"9.524858937833452E-08"
$wsOutput.Cells.Item(16, 5).Value = [double]"221866.8244862198"
$wsOutput.Cells.Item(16, 9).Value = [double]"1.024816990789045E-09"
$wsOutput.Cells.Item(16, 13).Value = [double]"2.829889203408771E-26"
$wsOutput.Cells.Item(17, 3).Value = [double]"2.514926763253608E-07"
$wsOutput.Cells.Item(17, 6).Value = [double]"1.142632588766418E-09"
$wsOutput.Cells.Item(17, 9).Value = [double]"2.679183820413722E-05"
$wsOutput.Cells.Item(17, 10).Value = [double]"0.698051099982766"
$wsOutput.Cells.Item(17, 13).Value = [double]"7.398192492396857E-22"
$wsOutput.Cells.Item(18, 3).Value = [double]"2.831010845942313E-08"
$wsOutput.Cells.Item(18, 4).Value = [double]"8.392470679719371E-09"
$wsOutput.Cells.Item(18, 6).Value = [double]"4.796181668456203E-10"
$wsOutput.Cells.Item(18, 9).Value = [double]"3.015912258236822E-06"
$wsOutput.Cells.Item(18, 10).Value = [double]"0.1536599911639926"
$wsOutput.Cells.Item(18, 13).Value = [double]"8.328021114717624E-23"
$wsOutput.Cells.Item(19, 3).Value = [double]"3.346922977608896E-09"
$wsOutput.Cells.Item(19, 4).Value = [double]"3.968704202855601E-09"
$wsOutput.Cells.Item(19, 6).Value = [double]"4.862832753241384E-09"
$wsOutput.Cells.Item(19, 9).Value = [double]"3.565520086231008E-07"
$wsOutput.Cells.Item(19, 10).Value = [double]"0.105366846483339"
$wsOutput.Cells.Item(19, 13).Value = [double]"9.845686485733951E-24"
$wsOutput.Cells.Item(20, 3).Value = [double]"0.004355978850511016"
$wsOutput.Cells.Item(20, 4).Value = [double]"0.0439043566409581"
$wsOutput.Cells.Item(20, 6).Value = [double]"5.772246289397745"
$wsOutput.Cells.Item(20, 9).Value = [double]"0.464048028311376"
$wsOutput.Cells.Item(20, 13).Value = [double]"1.281403916001046E-17"
$wsOutput.Cells.Item(21, 3).Value = [double]"0.009619868774929919"
$wsOutput.Cells.Item(21, 4).Value = [double]"0.9524858938047793"
$wsOutput.Cells.Item(21, 6).Value = [double]"221866.8244863141"
$wsOutput.Cells.Item(21, 9).Value = [double]"1.024816990811768"
$wsOutput.Cells.Item(21, 13).Value = [double]"2.829889203471516E-17"

# --- Input_flows sheet updates ---
$wsInput.Cells.Item(2, 3).Value = [double]"5.047917983504061E-16"
$wsInput.Cells.Item(15, 3).Value = [double]"1.3219918090551E-32"
$wsInput.Cells.Item(20, 3).Value = [double]"1.321991806593662E-23"
$wsInput.Cells.Item(22, 5).Value = [double]"3.993020435635347E-13"
$wsInput.Cells.Item(23, 5).Value = [double]"8.448780326926443E-15"
$wsInput.Cells.Item(24, 5).Value = [double]"3.910117373860849E-15"
$wsInput.Cells.Item(25, 5).Value = [double]"1.950812408815889E-15"
$wsInput.Cells.Item(26, 5).Value = [double]"1.600773355679499E-12"
$wsInput.Cells.Item(27, 5).Value = [double]"1.190898042553938E-05"
$wsInput.Cells.Item(28, 5).Value = [double]"5.864139260990636E-06"
$wsInput.Cells.Item(29, 5).Value = [double]"3.908780615001784E-06"
$wsInput.Cells.Item(30, 5).Value = [double]"1.950812284656906E-06"
$wsInput.Cells.Item(31, 5).Value = [double]"0.001600773355712505"
$wsInput.Cells.Item(32, 5).Value = [double]"1.437264317458851E-13"
$wsInput.Cells.Item(33, 5).Value = [double]"1.086308834832752E-14"
$wsInput.Cells.Item(34, 5).Value = [double]"4.216832679318602E-16"
$wsInput.Cells.Item(35, 5).Value = [double]"4.640162704057854E-10"
$wsInput.Cells.Item(36, 5).Value = [double]"1.006135795753847E-09"
$wsInput.Cells.Item(37, 5).Value = [double]"2.619738364288883E-05"
$wsInput.Cells.Item(38, 5).Value = [double]"2.941050389773365E-06"
$wsInput.Cells.Item(39, 5).Value = [double]"3.479414049127687E-07"
$wsInput.Cells.Item(40, 5).Value = [double]"0.4640162843175683"
$wsInput.Cells.Item(41, 5).Value = [double]"1.006135795775437"
$wsInput.Cells.Item(42, 6).Value = [double]"7.587655650349048E-29"
$wsInput.Cells.Item(43, 6).Value = [double]"4.742878024268019E-30"
$wsInput.Cells.Item(44, 6).Value = [double]"1.876613738314637E-31"
$wsInput.Cells.Item(45, 6).Value = [double]"9.278696419097372E-32"
$wsInput.Cells.Item(46, 6).Value = [double]"1.358374793077706E-28"
$wsInput.Cells.Item(47, 6).Value = [double]"6.893295761571661E-22"
$wsInput.Cells.Item(48, 6).Value = [double]"2.788270668194614E-22"
$wsInput.Cells.Item(49, 6).Value = [double]"1.858617924477421E-22"
$wsInput.Cells.Item(50, 6).Value = [double]"9.280380274846091E-23"
$wsInput.Cells.Item(51, 6).Value = [double]"1.356979295809919E-19"
$wsInput.Cells.Item(52, 6).Value = [double]"1.932779664317545E-29"
$wsInput.Cells.Item(53, 6).Value = [double]"1.396518503927476E-30"
$wsInput.Cells.Item(54, 6).Value = [double]"2.061495658879831E-32"
$wsInput.Cells.Item(55, 6).Value = [double]"8.867460482723861E-27"
$wsInput.Cells.Item(56, 6).Value = [double]"1.957425482408151E-26"
$wsInput.Cells.Item(57, 6).Value = [double]"5.116955009159427E-22"
$wsInput.Cells.Item(58, 6).Value = [double]"5.759690402071028E-23"
$wsInput.Cells.Item(59, 6).Value = [double]"6.809428120955666E-24"
$wsInput.Cells.Item(60, 6).Value = [double]"8.86773071530629E-18"
$wsInput.Cells.Item(61, 6).Value = [double]"1.957485058760349E-17"
